$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (@@ -932,22 +932,22 @@)
$ws.Range("H6").Value = 5717.778
$ws.Range("I6").Value = 6401.25
$ws.Range("K6").Value = 19203.75
$ws.Range("M6").Value = -19091.75

# Row 9 (@@ -1082,25 +1082,25 @@)
$ws.Range("H9").Value = 231
$ws.Range("I9").Value = 337.5
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 337.5
$ws.Range("L9").Value = 160
$ws.Range("M9").Value = -168.5
$ws.Range("N9").Value = -498

# Row 12 (@@ -1229,22 +1229,22 @@)
$ws.Range("H12").Value = 225
$ws.Range("I12").Value = 225
$ws.Range("K12").Value = 225
$ws.Range("M12").Value = -55

# Row 21 (@@ -1661,20 +1661,23 @@)
$ws.Range("H21").Value = 23673
$ws.Range("I21").Value = 27000
$ws.Range("K21").Value = 27000
$ws.Range("M21").Value = -26532

# Row 23 (@@ -1759,20 +1762,23 @@)
$ws.Range("H23").Value = 23673
$ws.Range("I23").Value = 27000
$ws.Range("K23").Value = 27000
$ws.Range("M23").Value = -26766

# Row 29 (@@ -2044,25 +2050,25 @@)
$ws.Range("H29").Value = 1746.4375
$ws.Range("J29").Value = 2095.6155
$ws.Range("L29").Value = 6286.8465
$ws.Range("N29").Value = -6848.8465

# Row 38 (@@ -2491,25 +2497,25 @@)
$ws.Range("H38").Value = 1849.6666
$ws.Range("I38").Value = 304
$ws.Range("J38").Value = 2622.5
$ws.Range("K38").Value = 912
$ws.Range("L38").Value = 7867.5
$ws.Range("M38").Value = -540
$ws.Range("N38").Value = -8611.5

# Row 58 (@@ -3483,25 +3489,25 @@)
$ws.Range("H58").Value = 2413.9333
$ws.Range("I58").Value = 524
$ws.Range("J58").Value = 4573.857
$ws.Range("K58").Value = 1572
$ws.Range("L58").Value = 13721.571
$ws.Range("M58").Value = -1422
$ws.Range("N58").Value = -14021.571

# Row 62 (@@ -3682,25 +3688,25 @@)
$ws.Range("H62").Value = 12350568
$ws.Range("I62").Value = 15877302
$ws.Range("J62").Value = 6999.5
$ws.Range("K62").Value = 15877302
$ws.Range("L62").Value = 6999.5
$ws.Range("M62").Value = -15876678
$ws.Range("N62").Value = -8247.5

# Row 65 (@@ -3832,25 +3838,25 @@)
$ws.Range("H65").Value = 12350568
$ws.Range("I65").Value = 15877302
$ws.Range("J65").Value = 6999.5
$ws.Range("K65").Value = 79386510
$ws.Range("L65").Value = 34997.5
$ws.Range("M65").Value = -79383390
$ws.Range("N65").Value = -41237.5

# Row 87 (@@ -4928,22 +4934,22 @@)
$ws.Range("H87").Value = 40863.285
$ws.Range("J87").Value = 40863.285
$ws.Range("L87").Value = 40863.285
$ws.Range("N87").Value = -43359.285

# Row 90 (@@ -5081,22 +5087,22 @@)
$ws.Range("H90").Value = 40863.285
$ws.Range("J90").Value = 40863.285
$ws.Range("L90").Value = 122589.855
$ws.Range("N90").Value = -135069.855

# Row 98 (@@ -5479,25 +5485,25 @@)
$ws.Range("H98").Value = 2909.375
$ws.Range("J98").Value = 988
$ws.Range("L98").Value = 988
$ws.Range("N98").Value = -3984

# Row 122 (@@ -6688,25 +6694,25 @@)
$ws.Range("H122").Value = 2909.375
$ws.Range("J122").Value = 988
$ws.Range("L122").Value = 2964
$ws.Range("N122").Value = -7864

# Row 135 (@@ -7334,22 +7340,22 @@)
$ws.Range("H135").Value = 27027452
$ws.Range("I135").Value = 230.0303
$ws.Range("K135").Value = 2070.2727
$ws.Range("M135").Value = 464.7273

# Row 137 (@@ -7435,22 +7441,22 @@)
$ws.Range("H137").Value = 1072.4247
$ws.Range("I137").Value = 816.8125
$ws.Range("K137").Value = 2450.4375
$ws.Range("M137").Value = 99.5625

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (@@ -9235,22 +9241,22 @@)
$ws.Range("H32").Value = 5104.8643
$ws.Range("I32").Value = 4701.4707
$ws.Range("K32").Value = 4701.4707
$ws.Range("M32").Value = -4414.4707

# Row 61 (@@ -10638,25 +10644,25 @@)
$ws.Range("H61").Value = 33334768
$ws.Range("I61").Value = 41667940
$ws.Range("J61").Value = 2085.6667
$ws.Range("K61").Value = 41667940
$ws.Range("L61").Value = 2085.6667
$ws.Range("M61").Value = -41667728
$ws.Range("N61").Value = -2509.6667

# Row 122 (@@ -13606,25 +13612,25 @@)
$ws.Range("H122").Value = 2397.9092
$ws.Range("I122").Value = 2890.5
$ws.Range("J122").Value = 1084.3334
$ws.Range("K122").Value = 8671.5
$ws.Range("L122").Value = 3253.0002
$ws.Range("M122").Value = -6221.5
$ws.Range("N122").Value = -8153.0002

# Row 132 (@@ -14090,22 +14096,22 @@)
$ws.Range("H132").Value = 2087.1516
$ws.Range("I132").Value = 2141.3845
$ws.Range("K132").Value = 6424.1535
$ws.Range("M132").Value = -3894.1535

# Row 136 (@@ -14283,25 +14289,25 @@)
$ws.Range("H136").Value = 33334768
$ws.Range("I136").Value = 41667940
$ws.Range("J136").Value = 2085.6667
$ws.Range("K136").Value = 125003820
$ws.Range("L136").Value = 6257.000100000001
$ws.Range("M136").Value = -125001270
$ws.Range("N136").Value = -11357.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 123 (@@ -20528,25 +20534,22 @@)
$ws.Range("H123").Value = 30486.666
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 30486.666
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30486.666
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -40286.666

# Row 134 (@@ -21064,22 +21067,22 @@)
$ws.Range("H134").Value = 4150.4443
$ws.Range("I134").Value = 1076.2667
$ws.Range("K134").Value = 3228.800099999999
$ws.Range("M134").Value = -693.8000999999995

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (@@ -22986,22 +22989,22 @@)
$ws.Range("H31").Value = 2220.4583
$ws.Range("I31").Value = 2294.2632
$ws.Range("K31").Value = 2294.2632
$ws.Range("M31").Value = -1999.2632

# Row 34 (@@ -23130,22 +23133,22 @@)
$ws.Range("H34").Value = 2220.4583
$ws.Range("I34").Value = 2294.2632
$ws.Range("K34").Value = 2294.2632
$ws.Range("M34").Value = -2092.2632

# Row 122 (@@ -27385,25 +27388,25 @@)
$ws.Range("H122").Value = 1125
$ws.Range("I122").Value = 1115.4546
$ws.Range("J122").Value = 1160
$ws.Range("K122").Value = 3346.3638
$ws.Range("L122").Value = 3480
$ws.Range("M122").Value = -896.3638000000001
$ws.Range("N122").Value = -8380

$ws = $wb.Worksheets.Item("CUL")
# Row 39 (@@ -30338,25 +30341,25 @@)
$ws.Range("H39").Value = 2566.2693
$ws.Range("J39").Value = 2469.2273
$ws.Range("L39").Value = 7407.6819
$ws.Range("N39").Value = -7995.6819

# Row 122 (@@ -34558,25 +34561,25 @@)
$ws.Range("H122").Value = 789.03845
$ws.Range("I122").Value = 613.6
$ws.Range("J122").Value = 898.6875
$ws.Range("K122").Value = 5522.400000000001
$ws.Range("L122").Value = 8088.1875
$ws.Range("M122").Value = -3072.400000000001
$ws.Range("N122").Value = -12988.1875

# Row 131 (@@ -35020,25 +35023,25 @@)
$ws.Range("H131").Value = 22728594
$ws.Range("J131").Value = 1566.8823
$ws.Range("L131").Value = 4700.6469
$ws.Range("N131").Value = -14780.6469

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (@@ -41690,22 +41693,22 @@)
$ws.Range("H126").Value = 2425
$ws.Range("I126").Value = 1850
$ws.Range("K126").Value = 5550
$ws.Range("M126").Value = -3080

# Row 132 (@@ -41984,22 +41987,22 @@)
$ws.Range("H132").Value = 2260.7036
$ws.Range("I132").Value = 1884.0454
$ws.Range("K132").Value = 5652.1362
$ws.Range("M132").Value = -3122.1362

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (@@ -44391,22 +44394,22 @@)
$ws.Range("H40").Value = 3334.3572
$ws.Range("I40").Value = 2052.3845
$ws.Range("K40").Value = 2052.3845
$ws.Range("M40").Value = -1916.3845

# Row 93 (@@ -46943,25 +46946,22 @@)
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 248
$ws.Range("N93").ClearContents()

# Row 122 (@@ -48337,25 +48337,25 @@)
$ws.Range("H122").Value = 22741618
$ws.Range("J122").Value = 8683.333000000001
$ws.Range("L122").Value = 26049.999
$ws.Range("N122").Value = -30949.999

# Row 132 (@@ -48830,22 +48830,22 @@)
$ws.Range("H132").Value = 28195.053
$ws.Range("I132").Value = 1675.3334
$ws.Range("K132").Value = 5026.0002
$ws.Range("M132").Value = -2496.0002

# Row 136 (@@ -49029,25 +49029,25 @@)
$ws.Range("H136").Value = 2616.082
$ws.Range("I136").Value = 2478.8965
$ws.Range("J136").Value = 5268.3335
$ws.Range("K136").Value = 7436.689499999999
$ws.Range("L136").Value = 15805.0005
$ws.Range("M136").Value = -4886.689499999999
$ws.Range("N136").Value = -20905.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 5 (@@ -49570,22 +49570,22 @@)
$ws.Range("H5").Value = 43000000
$ws.Range("J5").Value = 43000000
$ws.Range("L5").Value = 43000000
$ws.Range("N5").Value = -43000224

# Row 122 (@@ -55225,25 +55225,25 @@)
$ws.Range("H122").Value = 15629834
$ws.Range("I122").Value = 19236434
$ws.Range("J122").Value = 1235
$ws.Range("K122").Value = 57709302
$ws.Range("L122").Value = 3705
$ws.Range("M122").Value = -57706852
$ws.Range("N122").Value = -8605

# Row 126 (@@ -55424,22 +55424,22 @@)
$ws.Range("H126").Value = 40001532
$ws.Range("I126").Value = 55556844
$ws.Range("K126").Value = 166670532
$ws.Range("M126").Value = -166668062

# Row 132 (@@ -55718,25 +55718,25 @@)
$ws.Range("H132").Value = 6556.385
$ws.Range("I132").Value = 10954.714
$ws.Range("J132").Value = 1425
$ws.Range("K132").Value = 32864.142
$ws.Range("L132").Value = 4275
$ws.Range("M132").Value = -30334.142
$ws.Range("N132").Value = -9335

# Row 136 (@@ -55917,25 +55917,25 @@)
$ws.Range("H136").Value = 423.88235
$ws.Range("I136").Value = 341.6207
$ws.Range("J136").Value = 901
$ws.Range("K136").Value = 1024.8621
$ws.Range("L136").Value = 2703
$ws.Range("M136").Value = 1525.1379
$ws.Range("N136").Value = -7803
